$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 303-402 correspond to years 1901-2000: scale by 0.8
for ($r = 303; $r -le 402; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 * 0.8
}

# Rows 403-452 correspond to years 2001-2050: scale by 1.21
for ($r = 403; $r -le 452; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 * 1.21
}
